$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H93").Value = 41750
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 41750
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 41750
$ws.Range("N93").Value = -46742
$ws.Range("H95").Value = 19984
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 19984
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 19984
$ws.Range("N95").Value = -25476
$ws.Range("H120").Value = 35300
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 35300
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 35300
$ws.Range("N120").Value = -44976
$ws.Range("H134").Value = 125071410
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 125071410
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 125071410
$ws.Range("N134").Value = -125081550
$ws.Range("H135").Value = 886.89746
$ws.Range("I135").Value = 845.6129
$ws.Range("J135").Value = 1046.875
$ws.Range("K135").Value = 7610.5161
$ws.Range("L135").Value = 9421.875
$ws.Range("M135").Value = -5075.5161
$ws.Range("N135").Value = -14491.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 4003334.2
$ws.Range("I11").Value = 6000001.5
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 6000001.5
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -5999857.5
$ws.Range("N11").Value = -10288
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H95").Value = 25748.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 25748.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 25748.5
$ws.Range("N95").Value = -31240.5
$ws.Range("H96").Value = 12855.3
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 12855.3
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 12855.3
$ws.Range("N96").Value = -18347.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 18333634
$ws.Range("I6").Value = 36666830
$ws.Range("J6").Value = 434.66666
$ws.Range("K6").Value = 36666830
$ws.Range("L6").Value = 434.66666
$ws.Range("M6").Value = -36666717
$ws.Range("N6").Value = -660.66666
$ws.Range("H7").Value = 7215.5713
$ws.Range("I7").Value = 10049
$ws.Range("J7").Value = 132
$ws.Range("K7").Value = 10049
$ws.Range("L7").Value = 132
$ws.Range("M7").Value = -9936
$ws.Range("N7").Value = -358
$ws.Range("H15").Value = 9000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 9000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 9000
$ws.Range("N15").Value = -9340
$ws.Range("H17").Value = 2750
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2750
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2750
$ws.Range("N17").Value = -3098
$ws.Range("H19").Value = 629.8570999999999
$ws.Range("I19").Value = 629.8570999999999
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 629.8570999999999
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -459.8570999999999
$ws.Range("H24").Value = 629.8570999999999
$ws.Range("I24").Value = 629.8570999999999
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 629.8570999999999
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -459.8570999999999
$ws.Range("H25").Value = 70509.75
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 70509.75
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 70509.75
$ws.Range("N25").Value = -70857.75
$ws.Range("H41").Value = 999.3333
$ws.Range("I41").Value = 999.3333
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 999.3333
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -571.3333
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 30728.428
$ws.Range("I51").Value = 40000
$ws.Range("J51").Value = 29183.166
$ws.Range("K51").Value = 40000
$ws.Range("L51").Value = 29183.166
$ws.Range("M51").Value = -39264
$ws.Range("N51").Value = -30655.166
$ws.Range("H58").Value = 1067943.9
$ws.Range("I58").Value = 2955.375
$ws.Range("J58").Value = 2179236.2
$ws.Range("K58").Value = 2955.375
$ws.Range("L58").Value = 2179236.2
$ws.Range("M58").Value = -2752.375
$ws.Range("N58").Value = -2179642.2
$ws.Range("H59").Value = 34875
$ws.Range("I59").Value = 10000
$ws.Range("J59").Value = 43166.668
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 43166.668
$ws.Range("M59").Value = -8855
$ws.Range("N59").Value = -45456.668
$ws.Range("H60").Value = 16000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 16000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 16000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -17022
$ws.Range("H61").Value = 30728.428
$ws.Range("I61").Value = 40000
$ws.Range("J61").Value = 29183.166
$ws.Range("K61").Value = 40000
$ws.Range("L61").Value = 29183.166
$ws.Range("M61").Value = -39652
$ws.Range("N61").Value = -29879.166
$ws.Range("H95").Value = 7074.4443
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 7074.4443
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 7074.4443
$ws.Range("N95").Value = -12566.4443
$ws.Range("H96").Value = 23040.445
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 23040.445
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 23040.445
$ws.Range("N96").Value = -28532.445
$ws.Range("H136").Value = 1067943.9
$ws.Range("I136").Value = 2955.375
$ws.Range("J136").Value = 2179236.2
$ws.Range("K136").Value = 8866.125
$ws.Range("L136").Value = 6537708.600000001
$ws.Range("M136").Value = -6316.125
$ws.Range("N136").Value = -6542808.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3700
$ws.Range("I76").Value = 3400
$ws.Range("J76").Value = 3750
$ws.Range("K76").Value = 10200
$ws.Range("L76").Value = 11250
$ws.Range("M76").Value = -9817
$ws.Range("N76").Value = -12016
$ws.Range("H79").Value = 3700
$ws.Range("I79").Value = 3400
$ws.Range("J79").Value = 3750
$ws.Range("K79").Value = 10200
$ws.Range("L79").Value = 11250
$ws.Range("M79").Value = -8874
$ws.Range("N79").Value = -13902
$ws.Range("H126").Value = 2318.3333
$ws.Range("I126").Value = 965
$ws.Range("J126").Value = 2995
$ws.Range("K126").Value = 2895
$ws.Range("L126").Value = 8985
$ws.Range("M126").Value = 2045
$ws.Range("N126").Value = -18865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H95").Value = 9112.714
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 9112.714
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 9112.714
$ws.Range("N95").Value = -14604.714
$ws.Range("H140").Value = 28694.375
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 43185
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 43185
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -53545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 17387
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 17387
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 17387
$ws.Range("N95").Value = -22879
$ws.Range("H97").Value = 11654.571
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 11654.571
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 11654.571
$ws.Range("N97").Value = -13636.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 70871.60000000001
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 70871.60000000001
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 70871.60000000001
$ws.Range("N46").Value = -71333.60000000001
$ws.Range("H69").Value = 6855.2666
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 6855.2666
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 6855.2666
$ws.Range("N69").Value = -8353.266599999999
$ws.Range("H72").Value = 6855.2666
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 6855.2666
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 20565.7998
$ws.Range("N72").Value = -28053.7998
$ws.Range("H80").Value = 16433.666
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 16433.666
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 16433.666
$ws.Range("N80").Value = -18429.666
$ws.Range("H82").Value = 17000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 17000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 17000
$ws.Range("N82").Value = -17766
$ws.Range("H83").Value = 16433.666
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 16433.666
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 49300.99800000001
$ws.Range("N83").Value = -59284.99800000001
$ws.Range("H85").Value = 17000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 17000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 17000
$ws.Range("N85").Value = -19652
$ws.Range("H94").Value = 23443.334
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 23443.334
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 23443.334
$ws.Range("N94").Value = -25245.334
$ws.Range("H134").Value = 70871.60000000001
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 70871.60000000001
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 212614.8
$ws.Range("N134").Value = -217684.8
$ws.Range("H136").Value = 571358.7
$ws.Range("I136").Value = 756945.25
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2270835.75
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2268285.75
$ws.Range("N136").Value = -14100
